$wb = $excel.ActiveWorkbook

# Week 16 simulation added a new TE, J.Horsted, on the Bears roster.
$ws = $wb.Worksheets.Item("TE")

# Row 6: name in column A, simulated stat line (all zero so far) in B:J.
$ws.Cells.Item(6, 1).Value = "J.Horsted"
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 0

# Leave the TE tab active/selected with J7 highlighted, as the author left it.
$ws.Activate()
$ws.Range("J7").Select()
